$wb = $excel.ActiveWorkbook

# Insert the new "FTT-IH-CHI" worksheet right after "FTT-H" (before "FTT-S")
$ftH = $wb.Worksheets.Item("FTT-H")
$ws = $wb.Worksheets.Add($null, $ftH)
$ws.Name = "FTT-IH-CHI"

$rows = @(
    ,@("Variable name", "Read in?", "Code", "Description", "RowDim", "ColDim", "3DDim", "Conversion?", "Scenario")
    ,@("IUD1", 1, 0, "FTT-IH-CHI useful energy demand (GWh)", "ITTI", "TIME", "RSHORTTI", 0, "S0")
    ,@("ISC1", 1, 0, "FTT-IH-CHI market share caps", "ITTI", "TIME", "RSHORTTI", 0, "S0")
    ,@("IWW1", 1, 0, "FTT-IH-CHI cumulative capacities (EU28)", "ITTI", 0, 0, "TIME", "S0")
    ,@("IWI1", 1, 0, "FTT-IH-CHI yearly capacity additions", "ITTI", "TIME", "RSHORTTI", 0, "S0")
    ,@("IWK1", 1, 0, "FTT-IH-CHI yearly capacity ", "ITTI", "TIME", "RSHORTTI", 0, "S0")
    ,@("IWA1", 1, 0, "FTT-IH-CHI substitution matrix", "ITTI", "ITTI", "NA", 0, "S0")
    ,@("BIC1", 1, 0, "FTT-IH-CHI cost matrix", "ITTI", "CTTI", "RSHORTTI", 0, "S0")
    ,@("IWS1", 1, 0, "FTT-IH-CHI market shares", "ITTI", "TIME", "RSHORTTI", 0, "S0")
    ,@("IWB1", 1, 0, "FTT-IH-CHI learning spillover matrix", "ITTI", "ITTI", "NA", 0, "S0")
    ,@("IFD1", 1, 0, "FTT-IH-CHI final energy demand", "ITTI", "TIME", "RSHORTTI", 0, "S0")
    ,@("ILC1", 1, 0, "FTT-IH-CHI The real bare LC without taxes", "ITTI", "TIME", "RSHORTTI", 0, "S0")
    ,@("ILG1", 1, 0, "FTT-IH-CHI LC as seen by consumer", "ITTI", "TIME", "RSHORTTI", 0, "S0")
    ,@("ILD1", 1, 0, "FTT-IH-CHI LC standard deviation", "ITTI", "TIME", "RSHORTTI", 0, "S0")
    ,@("IWE1", 1, 0, "FTT-IH-CHI Emissions", "ITTI", "TIME", "RSHORTTI", 0, "S0")
    ,@("IHW1", 1, 0, "FTT-IH-CHI Global average emissions per UED (kt of CO2/GWh)", "ITTI", "TIME", "NA", 0, "S0")
    ,@("IAM1", 1, 0, "FTT-IH-CHI gamma values", "ITTI", "TIME", "RSHORTTI", 0, "S0")
    ,@("IRG1", 1, 0, "FTT-IH-CHI regulations (based on capacity)", "ITTI", "TIME", "RSHORTTI", 0, "S0")
    ,@("ISB1", 1, 0, "FTT-IH-CHI subsidies (percentage of investment cost)", "ITTI", "TIME", "RSHORTTI", 0, "S0")
    ,@("IXS1", 1, 0, "FTT-IH-CHI exogenous share changes", "ITTI", "TIME", "RSHORTTI", 0, "S0")
    ,@("IHF1", 1, 0, "FTT-IH-CHI final fuel demand for industrial heat (ktoe)", "JTI", "TIME", "RSHORTTI", 0, "S0")
    ,@("IJT1", 1, 0, "FTT-IH-CHI Tech to fuel conversion matrix (fuel x technology)", "JTI", "ITTI", "NA", 0, "S0")
    ,@("IFT1", 1, 0, "FTT-IH-CHI Fuel tax (2010 Euros/MWh) ", "ITTI", "TIME", "RSHORTTI", 0, "S0")
)


for ($r = 0; $r -lt $rows.Count; $r++) {
    $rowVals = $rows[$r]
    for ($c = 0; $c -lt $rowVals.Count; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $rowVals[$c]
    }
}

# Match the recorded selection/active cell on the new sheet
$ws.Range("J16").Select()

# Make the new sheet the active tab
$ws.Activate()

